$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns C:J are constant across every job row - read them once from row 1
$constVals = @()
for ($c = 3; $c -le 10; $c++) {
    $constVals += $ws.Cells.Item(1, $c).Value2
}

# Column A alternates between these two image file names
$imgA = "blog_1.jpg"
$imgB = "blog_4.jpg"

# Update existing rows 2-4 (Column B changes) and add new rows 5-11
$jobRows = @(
    @{ Row = 2;  A = $imgB; B = "Sales Officer" },
    @{ Row = 3;  A = $imgA; B = "Field Officer" },
    @{ Row = 4;  A = $imgB; B = "Marketing" },
    @{ Row = 5;  A = $imgA; B = "Production, Soldering" },
    @{ Row = 6;  A = $imgB; B = "Sales Officer" },
    @{ Row = 7;  A = $imgA; B = "Field Officer" },
    @{ Row = 8;  A = $imgB; B = "Marketing" },
    @{ Row = 9;  A = $imgB; B = "Sales Officer" },
    @{ Row = 10; A = $imgA; B = "Field Officer" },
    @{ Row = 11; A = $imgB; B = "Marketing" }
)

foreach ($jr in $jobRows) {
    $r = $jr.Row
    $ws.Cells.Item($r, 1).Value = $jr.A
    $ws.Cells.Item($r, 2).Value = $jr.B
    for ($c = 3; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $constVals[$c - 3]
    }
}

# Match the final selection recorded in the saved workbook
$ws.Range("A9:XFD11").Select()
